# Applies the "Penality Reward System" data edit:
#   - On the "Weekly Quantity" sheet, two obsolete weekly rows (originally
#     rows 11 and 12) are removed, shifting all subsequent rows up by two
#     and shrinking the used range from A1:B49 to A1:B47.
#   - On the "Monthly Trend" sheet, the requested quantity for the month
#     in row 4 is corrected from 320 to 190.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows.Item(11).Delete()
$wsWeekly.Rows.Item(11).Delete()

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B4").Value = 190
